$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add hyperlink on K3 pointing to the ACM To BTK converter repository,
# with the URL text itself shown in the cell (matches target shared string).
$ws.Hyperlinks.Add($ws.Range("K3"), "https://github.com/btk42/CV-ATB-00-000-MCR-S2B0-01")

# Update the active selection to F10, as left by the author after editing.
$ws.Range("F10").Select()
